$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F ("想去人数" / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 179
$ws1.Range("F5").Value = 281
$ws1.Range("F6").Value = 391
$ws1.Range("F7").Value = 245
$ws1.Range("F8").Value = 2315
$ws1.Range("F10").Value = 5752
$ws1.Range("F12").Value = 377

# Sheet "全部类型" (All Types) - update column F ("想去人数" / interested count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 179
$ws4.Range("F6").Value = 281
$ws4.Range("F7").Value = 391
$ws4.Range("F8").Value = 245
$ws4.Range("F11").Value = 2315
$ws4.Range("F13").Value = 5752
$ws4.Range("F15").Value = 377
